# === Re-plot toy-spam confidence anchors with "min 5" occurrence cutoff ===
# Every anchor-word statistic is recomputed against the new run; the
# "negative" table also loses its trailing row (old row 31, "like").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the row that no longer exists in the re-plotted output
$ws.Rows.Item(31).Delete() | Out-Null

# --- "negative" table (name/anchor score/.../normal => columns A:H) ---
$negArr = New-Object "object[,]" 30,8
$negArr[0,0] = "negative"
$negArr[1,0] = "name"
$negArr[1,1] = "anchor score"
$negArr[1,2] = "type occurences"
$negArr[1,3] = "total occurences"
$negArr[1,4] = "+%"
$negArr[1,5] = "-%"
$negArr[1,6] = "both"
$negArr[1,7] = "normal"
$negArr[2,0] = "poorly"
$negArr[2,1] = 1
$negArr[2,2] = 46
$negArr[2,3] = 46
$negArr[2,4] = 0
$negArr[2,5] = 1
$negArr[2,6] = $false
$negArr[2,7] = 0
$negArr[3,0] = "however"
$negArr[3,1] = 0.78125
$negArr[3,2] = 50
$negArr[3,3] = 50
$negArr[3,4] = 0
$negArr[3,5] = 1
$negArr[3,6] = $false
$negArr[3,7] = 14
$negArr[4,0] = "disappointing"
$negArr[4,1] = 0.7727272727272727
$negArr[4,2] = 34
$negArr[4,3] = 34
$negArr[4,4] = 0
$negArr[4,5] = 1
$negArr[4,6] = $false
$negArr[4,7] = 10
$negArr[5,0] = "poor"
$negArr[5,1] = 0.7323943661971831
$negArr[5,2] = 52
$negArr[5,3] = 52
$negArr[5,4] = 0
$negArr[5,5] = 1
$negArr[5,6] = $false
$negArr[5,7] = 19
$negArr[6,0] = "disappointed"
$negArr[6,1] = 0.7258064516129032
$negArr[6,2] = 135
$negArr[6,3] = 135
$negArr[6,4] = 0
$negArr[6,5] = 1
$negArr[6,6] = $false
$negArr[6,7] = 51
$negArr[7,0] = "broke"
$negArr[7,1] = 0.7087378640776699
$negArr[7,2] = 146
$negArr[7,3] = 146
$negArr[7,4] = 0
$negArr[7,5] = 1
$negArr[7,6] = $false
$negArr[7,7] = 60
$negArr[8,0] = "instead"
$negArr[8,1] = 0.625
$negArr[8,2] = 30
$negArr[8,3] = 30
$negArr[8,4] = 0
$negArr[8,5] = 1
$negArr[8,6] = $false
$negArr[8,7] = 18
$negArr[9,0] = "waste"
$negArr[9,1] = 0.6148648648648649
$negArr[9,2] = 91
$negArr[9,3] = 91
$negArr[9,4] = 0
$negArr[9,5] = 1
$negArr[9,6] = $false
$negArr[9,7] = 57
$negArr[10,0] = "smaller"
$negArr[10,1] = 0.6050420168067226
$negArr[10,2] = 72
$negArr[10,3] = 72
$negArr[10,4] = 0
$negArr[10,5] = 1
$negArr[10,6] = $false
$negArr[10,7] = 47
$negArr[11,0] = "junk"
$negArr[11,1] = 0.5636363636363636
$negArr[11,2] = 31
$negArr[11,3] = 31
$negArr[11,4] = 0
$negArr[11,5] = 1
$negArr[11,6] = $false
$negArr[11,7] = 24
$negArr[12,0] = "small"
$negArr[12,1] = 0.4956521739130435
$negArr[12,2] = 171
$negArr[12,3] = 171
$negArr[12,4] = 0
$negArr[12,5] = 1
$negArr[12,6] = $false
$negArr[12,7] = 174
$negArr[13,0] = "apart"
$negArr[13,1] = 0.4526315789473684
$negArr[13,2] = 43
$negArr[13,3] = 43
$negArr[13,4] = 0
$negArr[13,5] = 1
$negArr[13,6] = $false
$negArr[13,7] = 52
$negArr[14,0] = "plastic"
$negArr[14,1] = 0.4173228346456693
$negArr[14,2] = 53
$negArr[14,3] = 53
$negArr[14,4] = 0
$negArr[14,5] = 1
$negArr[14,6] = $false
$negArr[14,7] = 74
$negArr[15,0] = "broken"
$negArr[15,1] = 0.4096385542168675
$negArr[15,2] = 34
$negArr[15,3] = 34
$negArr[15,4] = 0
$negArr[15,5] = 1
$negArr[15,6] = $false
$negArr[15,7] = 49
$negArr[16,0] = "thought"
$negArr[16,1] = 0.297029702970297
$negArr[16,2] = 60
$negArr[16,3] = 60
$negArr[16,4] = 0
$negArr[16,5] = 1
$negArr[16,6] = $false
$negArr[16,7] = 142
$negArr[17,0] = "ok"
$negArr[17,1] = 0.28125
$negArr[17,2] = 36
$negArr[17,3] = 36
$negArr[17,4] = 0
$negArr[17,5] = 1
$negArr[17,6] = $false
$negArr[17,7] = 92
$negArr[18,0] = "cheap"
$negArr[18,1] = 0.2748815165876777
$negArr[18,2] = 58
$negArr[18,3] = 58
$negArr[18,4] = 0
$negArr[18,5] = 1
$negArr[18,6] = $false
$negArr[18,7] = 153
$negArr[19,0] = "size"
$negArr[19,1] = 0.2319587628865979
$negArr[19,2] = 45
$negArr[19,3] = 45
$negArr[19,4] = 0
$negArr[19,5] = 1
$negArr[19,6] = $false
$negArr[19,7] = 149
$negArr[20,0] = "hard"
$negArr[20,1] = 0.185
$negArr[20,2] = 37
$negArr[20,3] = 37
$negArr[20,4] = 0
$negArr[20,5] = 1
$negArr[20,6] = $false
$negArr[20,7] = 163
$negArr[21,0] = "item"
$negArr[21,1] = 0.1811594202898551
$negArr[21,2] = 50
$negArr[21,3] = 50
$negArr[21,4] = 0
$negArr[21,5] = 1
$negArr[21,6] = $false
$negArr[21,7] = 226
$negArr[22,0] = "money"
$negArr[22,1] = 0.1708860759493671
$negArr[22,2] = 54
$negArr[22,3] = 54
$negArr[22,4] = 0
$negArr[22,5] = 1
$negArr[22,6] = $false
$negArr[22,7] = 262
$negArr[23,0] = "would"
$negArr[23,1] = 0.1679049034175334
$negArr[23,2] = 113
$negArr[23,3] = 114
$negArr[23,4] = 0.01
$negArr[23,5] = 0.99
$negArr[23,6] = $true
$negArr[23,7] = 560
$negArr[24,0] = "used"
$negArr[24,1] = 0.1657142857142857
$negArr[24,2] = 29
$negArr[24,3] = 29
$negArr[24,4] = 0
$negArr[24,5] = 1
$negArr[24,6] = $false
$negArr[24,7] = 146
$negArr[25,0] = "work"
$negArr[25,1] = 0.1582278481012658
$negArr[25,2] = 50
$negArr[25,3] = 50
$negArr[25,4] = 0
$negArr[25,5] = 1
$negArr[25,6] = $false
$negArr[25,7] = 266
$negArr[26,0] = "product"
$negArr[26,1] = 0.13215859030837
$negArr[26,2] = 60
$negArr[26,3] = 60
$negArr[26,4] = 0
$negArr[26,5] = 1
$negArr[26,6] = $false
$negArr[26,7] = 394
$negArr[27,0] = "price"
$negArr[27,1] = 0.1206896551724138
$negArr[27,2] = 42
$negArr[27,3] = 42
$negArr[27,4] = 0
$negArr[27,5] = 1
$negArr[27,6] = $false
$negArr[27,7] = 306
$negArr[28,0] = "use"
$negArr[28,1] = 0.08767123287671233
$negArr[28,2] = 32
$negArr[28,3] = 32
$negArr[28,4] = 0
$negArr[28,5] = 1
$negArr[28,6] = $false
$negArr[28,7] = 333
$negArr[29,0] = "like"
$negArr[29,1] = 0.08552631578947369
$negArr[29,2] = 52
$negArr[29,3] = 52
$negArr[29,4] = 0
$negArr[29,5] = 1
$negArr[29,6] = $false
$negArr[29,7] = 556
$ws.Range("A1:H30").Value = $negArr

# --- "positive" table (name/anchor score/.../normal => columns J:Q) ---
$posArr = New-Object "object[,]" 30,8
$posArr[0,0] = "positive"
$posArr[1,0] = "name"
$posArr[1,1] = "anchor score"
$posArr[1,2] = "type occurences"
$posArr[1,3] = "total occurences"
$posArr[1,4] = "+%"
$posArr[1,5] = "-%"
$posArr[1,6] = "both"
$posArr[1,7] = "normal"
$posArr[2,0] = "wonderful"
$posArr[2,1] = 0.8392857142857143
$posArr[2,2] = 47
$posArr[2,3] = 47
$posArr[2,4] = 1
$posArr[2,5] = 0
$posArr[2,6] = $false
$posArr[2,7] = 9
$posArr[3,0] = "awesome"
$posArr[3,1] = 0.8153846153846154
$posArr[3,2] = 53
$posArr[3,3] = 53
$posArr[3,4] = 1
$posArr[3,5] = 0
$posArr[3,6] = $false
$posArr[3,7] = 12
$posArr[4,0] = "favorite"
$posArr[4,1] = 0.6451612903225806
$posArr[4,2] = 60
$posArr[4,3] = 60
$posArr[4,4] = 1
$posArr[4,5] = 0
$posArr[4,6] = $false
$posArr[4,7] = 33
$posArr[5,0] = "classic"
$posArr[5,1] = 0.5849056603773585
$posArr[5,2] = 31
$posArr[5,3] = 31
$posArr[5,4] = 1
$posArr[5,5] = 0
$posArr[5,6] = $false
$posArr[5,7] = 22
$posArr[6,0] = "excellent"
$posArr[6,1] = 0.484375
$posArr[6,2] = 31
$posArr[6,3] = 31
$posArr[6,4] = 1
$posArr[6,5] = 0
$posArr[6,6] = $false
$posArr[6,7] = 33
$posArr[7,0] = "great"
$posArr[7,1] = 0.3467213114754099
$posArr[7,2] = 423
$posArr[7,3] = 423
$posArr[7,4] = 1
$posArr[7,5] = 0
$posArr[7,6] = $false
$posArr[7,7] = 797
$posArr[8,0] = "love"
$posArr[8,1] = 0.3055954088952654
$posArr[8,2] = 213
$posArr[8,3] = 213
$posArr[8,4] = 1
$posArr[8,5] = 0
$posArr[8,6] = $false
$posArr[8,7] = 484
$posArr[9,0] = "loves"
$posArr[9,1] = 0.2427385892116183
$posArr[9,2] = 117
$posArr[9,3] = 117
$posArr[9,4] = 1
$posArr[9,5] = 0
$posArr[9,6] = $false
$posArr[9,7] = 365
$posArr[10,0] = "perfect"
$posArr[10,1] = 0.1987951807228916
$posArr[10,2] = 33
$posArr[10,3] = 33
$posArr[10,4] = 1
$posArr[10,5] = 0
$posArr[10,6] = $false
$posArr[10,7] = 133
$posArr[11,0] = "loved"
$posArr[11,1] = 0.1681957186544343
$posArr[11,2] = 55
$posArr[11,3] = 55
$posArr[11,4] = 1
$posArr[11,5] = 0
$posArr[11,6] = $false
$posArr[11,7] = 272
$posArr[12,0] = "fun"
$posArr[12,1] = 0.0736196319018405
$posArr[12,2] = 84
$posArr[12,3] = 84
$posArr[12,4] = 1
$posArr[12,5] = 0
$posArr[12,6] = $false
$posArr[12,7] = 1057
$posArr[13,0] = "game"
$posArr[13,1] = 0.04347826086956522
$posArr[13,2] = 67
$posArr[13,3] = 67
$posArr[13,4] = 1
$posArr[13,5] = 0
$posArr[13,6] = $false
$posArr[13,7] = 1474
$ws.Range("J1:Q30").Value = $posArr

Write-Host "Applied toy-spam min-5 data update"
